# Add the "2022-Q4" detail sheet (new quarter) right after "总计" and
# before the previously-first quarter sheet "2022-Q3"; update the "总计"
# (summary) sheet with the new quarter's aggregate row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet before "2022-Q3".
# ---------------------------------------------------------------------
$wsQ3Ref = $wb.Worksheets.Item("2022-Q3")
$newWs = $wb.Worksheets.Add($wsQ3Ref)
$newWs.Name = "2022-Q4"

# Re-fetch sheet references now that the collection changed - stale
# references left over from before an Add()/structural change lose their
# formatting info in this runtime.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Copy header-row (B1:H1) and first-data-row (A-column) formatting from
# the "2022-Q3" sheet so the new sheet matches the existing look (bold
# header font, centered/top aligned + thin-bordered index column, etc).
$wsQ3.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ3.Range("A2").Copy()
$newWs.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header labels.
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Data rows. Columns B-G are kept as literal text (fund codes keep
# leading zeros, the numeric-looking figures stay text) just like the
# source data; only A (row index) and H (rank) are real numbers.
$q4 = @(
    @(0, "005416", "鹏华尊惠18个月定期开放混合A",       "2.83", "39.19", "1.59", "0.0450", 9),
    @(1, "002123", "北信瑞丰外延增长主题灵活配置混合", "0.20", "93.59", "5.07", "0.0101", 5),
    @(2, "005417", "鹏华尊惠18个月定期开放混合C",       "0.36", "39.19", "1.59", "0.0057", 9),
    @(3, "001154", "北信瑞丰平安中国主题灵活配置混合", "0.13", "93.42", "3.99", "0.0052", 4)
)

for ($i = 0; $i -lt $q4.Length; $i++) {
    $r = $i + 2
    $row = $q4[$i]

    $newWs.Cells.Item($r, 1).Value = $row[0]

    $newWs.Cells.Item($r, 2).NumberFormat = "@"
    $newWs.Cells.Item($r, 2).Value = $row[1]

    $newWs.Cells.Item($r, 3).NumberFormat = "@"
    $newWs.Cells.Item($r, 3).Value = $row[2]

    $newWs.Cells.Item($r, 4).NumberFormat = "@"
    $newWs.Cells.Item($r, 4).Value = $row[3]

    $newWs.Cells.Item($r, 5).NumberFormat = "@"
    $newWs.Cells.Item($r, 5).Value = $row[4]

    $newWs.Cells.Item($r, 6).NumberFormat = "@"
    $newWs.Cells.Item($r, 6).Value = $row[5]

    $newWs.Cells.Item($r, 7).NumberFormat = "@"
    $newWs.Cells.Item($r, 7).Value = $row[6]

    $newWs.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: a new first data row for 2022-Q4,
#    with every following row's index (column A) bumped by one.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Grow the used range by one row, re-using the existing index-column
# (A) style for the freshly added row 10.
$wsTotal.Range("A9").Copy()
$wsTotal.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totals = @(
    @(0, "2022-Q4", 4,  0.07000000000000001),
    @(1, "2022-Q3", 11, 2.29),
    @(2, "2022-Q2", 9,  2.71),
    @(3, "2022-Q1", 11, 3.66),
    @(4, "2021-Q4", 20, 5.22),
    @(5, "2021-Q3", 6,  1.08),
    @(6, "2021-Q2", 11, 1.45),
    @(7, "2021-Q1", 6,  1.42),
    @(8, "2020-Q4", 7,  3.87)
)

for ($i = 0; $i -lt $totals.Length; $i++) {
    $r = $i + 2
    $row = $totals[$i]
    $wsTotal.Cells.Item($r, 1).Value = $row[0]
    $wsTotal.Cells.Item($r, 2).Value = $row[1]
    $wsTotal.Cells.Item($r, 3).Value = $row[2]
    $wsTotal.Cells.Item($r, 4).Value = $row[3]
}

# Restore "总计" as the active tab (adding the sheet above made the new
# "2022-Q4" sheet active instead), matching the workbook's original
# bookViews/activeTab state.
$wsTotal.Activate()
